$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting instead of
# being auto-converted to a number by Excel when values look numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "51.035.84"
$ws.Cells.Item(2, 5).Value = "  -0.67%  "
$ws.Cells.Item(3, 4).Value = "2.938.78"
$ws.Cells.Item(3, 5).Value = "  +0.67%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).Value = "379.30"
$ws.Cells.Item(5, 5).Value = "  +0.96%  "
$ws.Cells.Item(6, 4).Value = "101.85"
$ws.Cells.Item(6, 5).Value = "  -0.85%  "
$ws.Cells.Item(7, 4).Value = "0.537"
$ws.Cells.Item(7, 5).Value = "  -0.82%  "
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
$ws.Cells.Item(9, 4).Value = "0.585"
$ws.Cells.Item(9, 5).Value = "  +0.39%  "
$ws.Cells.Item(10, 4).Value = "36.18"
$ws.Cells.Item(10, 5).Value = "  -1.99%  "
$ws.Cells.Item(11, 5).Value = "  -0.19%  "
$ws.Cells.Item(12, 4).Value = "0.0835"
$ws.Cells.Item(12, 5).Value = "  +0.15%  "
$ws.Cells.Item(13, 4).Value = "3.411.74"
$ws.Cells.Item(13, 5).Value = "  +0.94%  "
$ws.Cells.Item(14, 4).Value = "17.96"
$ws.Cells.Item(14, 5).Value = "  -1.84%  "
$ws.Cells.Item(15, 4).Value = "7.41"
$ws.Cells.Item(15, 5).Value = "  +1.15%  "
$ws.Cells.Item(16, 4).Value = "2.940.70"
$ws.Cells.Item(16, 5).Value = "  +0.91%  "
$ws.Cells.Item(17, 4).Value = "0.979"
$ws.Cells.Item(17, 5).Value = "  +5.73%  "
$ws.Cells.Item(18, 4).Value = "50.927.46"
$ws.Cells.Item(18, 5).Value = "  -0.76%  "
$ws.Cells.Item(19, 4).Value = "3.21"
$ws.Cells.Item(19, 5).Value = "  -5.53%  "
$ws.Cells.Item(20, 4).Value = "7.28"
$ws.Cells.Item(20, 5).Value = "  -0.42%  "
$ws.Cells.Item(21, 4).Value = "12.48"
$ws.Cells.Item(21, 5).Value = "  -2.98%  "
$ws.Cells.Item(22, 4).Value = "0.0₃0950"
$ws.Cells.Item(22, 5).Value = "  +0.56%  "
$ws.Cells.Item(23, 4).Value = "68.35"
$ws.Cells.Item(23, 5).Value = "  +0.13%  "
$ws.Cells.Item(24, 4).Value = "260.27"
$ws.Cells.Item(24, 5).Value = "  -0.43%  "
$ws.Cells.Item(25, 4).Value = "2.86"
$ws.Cells.Item(25, 5).Value = "  +4.52%  "
$ws.Cells.Item(26, 4).Value = "8.06"
$ws.Cells.Item(26, 5).Value = "  +11.76%  "
$ws.Cells.Item(27, 4).Value = "7.48"
$ws.Cells.Item(27, 5).Value = "  +9.48%  "
$ws.Cells.Item(28, 4).Value = "4.10"
$ws.Cells.Item(28, 5).Value = "  -0.70%  "
$ws.Cells.Item(32, 4).Value = "25.57"
$ws.Cells.Item(32, 5).Value = "  -0.21%  "
$ws.Cells.Item(33, 4).Value = "9.74"
$ws.Cells.Item(33, 5).Value = "  -0.51%  "
$ws.Cells.Item(34, 4).Value = "50.49"
$ws.Cells.Item(34, 5).Value = "  -1.20%  "
$ws.Cells.Item(35, 4).Value = "33.80"
$ws.Cells.Item(35, 5).Value = "  -0.13%  "
$ws.Cells.Item(36, 4).Value = "2.04"
$ws.Cells.Item(36, 5).Value = "  -2.93%  "
$ws.Cells.Item(37, 4).Value = "0.0441"
$ws.Cells.Item(37, 5).Value = "  +4.81%  "
$ws.Cells.Item(38, 5).Value = "  -0.01%  "
$ws.Cells.Item(39, 4).Value = "2.96"
$ws.Cells.Item(39, 5).Value = "  -1.09%  "
$ws.Cells.Item(40, 4).Value = "16.85"
$ws.Cells.Item(40, 5).Value = "  -0.55%  "
$ws.Cells.Item(41, 4).Value = "2.53"
$ws.Cells.Item(41, 5).Value = "  -0.07%  "
$ws.Cells.Item(42, 5).Value = "  +0.85%  "
$ws.Cells.Item(43, 4).Value = "1.76"
$ws.Cells.Item(43, 5).Value = "  -2.47%  "
$ws.Cells.Item(44, 4).Value = "121.57"
$ws.Cells.Item(44, 5).Value = "  -0.75%  "
$ws.Cells.Item(45, 4).Value = "21.08"
$ws.Cells.Item(45, 5).Value = "  -1.95%  "
$ws.Cells.Item(46, 5).Value = "  +1.03%  "
$ws.Cells.Item(47, 4).Value = "0.272"
$ws.Cells.Item(47, 5).Value = "  +1.58%  "
$ws.Cells.Item(48, 5).Value = "  +1.83%  "
$ws.Cells.Item(49, 4).Value = "2.002.81"
$ws.Cells.Item(49, 5).Value = "  -0.78%  "
$ws.Cells.Item(50, 4).Value = "3.20"
$ws.Cells.Item(50, 5).Value = "  +1.71%  "
$ws.Cells.Item(51, 4).Value = "0.0333"
$ws.Cells.Item(51, 5).Value = "  +5.99%  "

# Rows 29-31 were reordered (Hedera/Dai/Kaspa -> Dai/Kaspa/Hedera) with updated data
$ws.Cells.Item(29, 2).Value = "Dai"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 5).Value = "  -0.01%  "
$ws.Cells.Item(30, 2).Value = "Kaspa"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(30, 4).Value = "0.165"
$ws.Cells.Item(30, 5).Value = "  -2.09%  "
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).Value = "0.112"
$ws.Cells.Item(31, 5).Value = "  +10.47%  "
